{"js": "// Replace the 25 \"a\u00f7b=c, d\" answer cells in the practice-sheet table.\n// The mapping below is applied in document order (top-to-bottom,\n// left-to-right through the 5 data rows) so the one duplicated source\n// value (\"42\u00f75=8, 2\") still resolves to the correct distinct targets.\nconst replacements = [\n  [\"86\u00f73=28, 2\", \"33\u00f76=5, 3\"],\n  [\"27\u00f75=5, 2\", \"26\u00f76=4, 2\"],\n  [\"83\u00f79=9, 2\", \"10\u00f76=1, 4\"],\n  [\"67\u00f75=13, 2\", \"99\u00f72=49, 1\"],\n  [\"55\u00f72=27, 1\", \"92\u00f79=10, 2\"],\n  [\"87\u00f76=14, 3\", \"64\u00f78=8, 0\"],\n  [\"90\u00f78=11, 2\", \"37\u00f76=6, 1\"],\n  [\"77\u00f79=8, 5\", \"99\u00f73=33, 0\"],\n  [\"86\u00f79=9, 5\", \"29\u00f77=4, 1\"],\n  [\"51\u00f77=7, 2\", \"95\u00f74=23, 3\"],\n  [\"28\u00f77=4, 0\", \"79\u00f74=19, 3\"],\n  [\"85\u00f74=21, 1\", \"99\u00f74=24, 3\"],\n  [\"37\u00f79=4, 1\", \"48\u00f79=5, 3\"],\n  [\"67\u00f79=7, 4\", \"11\u00f79=1, 2\"],\n  [\"22\u00f75=4, 2\", \"25\u00f79=2, 7\"],\n  [\"92\u00f76=15, 2\", \"42\u00f78=5, 2\"],\n  [\"42\u00f75=8, 2\", \"14\u00f76=2, 2\"],\n  [\"35\u00f75=7, 0\", \"11\u00f78=1, 3\"],\n  [\"30\u00f77=4, 2\", \"28\u00f75=5, 3\"],\n  [\"42\u00f75=8, 2\", \"46\u00f75=9, 1\"],\n  [\"59\u00f74=14, 3\", \"47\u00f78=5, 7\"],\n  [\"57\u00f76=9, 3\", \"22\u00f79=2, 4\"],\n  [\"56\u00f77=8, 0\", \"87\u00f79=9, 6\"],\n  [\"28\u00f75=5, 3\", \"21\u00f74=5, 1\"],\n  [\"34\u00f72=17, 0\", \"58\u00f74=14, 2\"],\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet cursor = 0;\nfor (const paragraph of paragraphs.items) {\n  if (cursor >= replacements.length) break;\n  const [oldText, newText] = replacements[cursor];\n  if (paragraph.text === oldText) {\n    paragraph.insertText(newText, \"Replace\");\n    cursor++;\n  }\n}\nawait context.sync();\n\nif (cursor !== replacements.length) {\n  throw new Error(\n    `Only matched ${cursor} of ${replacements.length} expected answer cells`\n  );\n}\n", "ps1": "# Replace the 25 \"a\u00f7b=c, d\" answer cells in the practice-sheet table.\n# The mapping below is applied in document order (top-to-bottom,\n# left-to-right through the 5 data rows) so the one duplicated source\n# value (\"42\u00f75=8, 2\") still resolves to the correct distinct targets.\n$oldValues = @(\n    \"86\u00f73=28, 2\",\n    \"27\u00f75=5, 2\",\n    \"83\u00f79=9, 2\",\n    \"67\u00f75=13, 2\",\n    \"55\u00f72=27, 1\",\n    \"87\u00f76=14, 3\",\n    \"90\u00f78=11, 2\",\n    \"77\u00f79=8, 5\",\n    \"86\u00f79=9, 5\",\n    \"51\u00f77=7, 2\",\n    \"28\u00f77=4, 0\",\n    \"85\u00f74=21, 1\",\n    \"37\u00f79=4, 1\",\n    \"67\u00f79=7, 4\",\n    \"22\u00f75=4, 2\",\n    \"92\u00f76=15, 2\",\n    \"42\u00f75=8, 2\",\n    \"35\u00f75=7, 0\",\n    \"30\u00f77=4, 2\",\n    \"42\u00f75=8, 2\",\n    \"59\u00f74=14, 3\",\n    \"57\u00f76=9, 3\",\n    \"56\u00f77=8, 0\",\n    \"28\u00f75=5, 3\",\n    \"34\u00f72=17, 0\"\n)\n$newValues = @(\n    \"33\u00f76=5, 3\",\n    \"26\u00f76=4, 2\",\n    \"10\u00f76=1, 4\",\n    \"99\u00f72=49, 1\",\n    \"92\u00f79=10, 2\",\n    \"64\u00f78=8, 0\",\n    \"37\u00f76=6, 1\",\n    \"99\u00f73=33, 0\",\n    \"29\u00f77=4, 1\",\n    \"95\u00f74=23, 3\",\n    \"79\u00f74=19, 3\",\n    \"99\u00f74=24, 3\",\n    \"48\u00f79=5, 3\",\n    \"11\u00f79=1, 2\",\n    \"25\u00f79=2, 7\",\n    \"42\u00f78=5, 2\",\n    \"14\u00f76=2, 2\",\n    \"11\u00f78=1, 3\",\n    \"28\u00f75=5, 3\",\n    \"46\u00f75=9, 1\",\n    \"47\u00f78=5, 7\",\n    \"22\u00f79=2, 4\",\n    \"87\u00f79=9, 6\",\n    \"21\u00f74=5, 1\",\n    \"58\u00f74=14, 2\"\n)\n\n$d = $word.ActiveDocument\n$cursor = 0\n\nforeach ($p in $d.Paragraphs) {\n    if ($cursor -ge $oldValues.Count) {\n        break\n    }\n    $r = $p.Range\n    $t = $r.Text\n    # Strip the trailing paragraph mark / cell mark so it compares equal to\n    # the plain answer text.\n    $t = $t.TrimEnd([char]13, [char]7)\n    if ($t -eq $oldValues[$cursor]) {\n        $r.Text = $newValues[$cursor]\n        $cursor++\n    }\n}\n\nif ($cursor -ne $oldValues.Count) {\n    throw \"Only matched $cursor of $($oldValues.Count) expected answer cells\"\n}\n\nWrite-Output \"Replaced $cursor answer cells\"\n"}
